$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Fin de viaje"
$ws.Range("B3").Value = "Virginia Woolf"
$ws.Range("C3").Value = "Lumen"
